$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Note" column: header in C1/D1, description text in D1
$ws.Range("C1").Value = "Note"
$ws.Range("D1").Value = 'expected xuong dong thi anh xu ly theo cach xuong dong (replace ''\n'' by ''${EMPTY}'')'

# Widen column D to match the author's new column width (~78.29 characters)
$ws.Columns.Item(4).ColumnWidth = 77.43

# Move the active selection to D7, matching the recorded cursor position
$ws.Range("D7").Select() | Out-Null
